    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Cells whose new value looks like a plain number (e.g. "108.05") need the
    # destination cell pre-formatted as Text, otherwise Excel's COM type-inference
    # would silently store them as numeric values and lose the original text
    # formatting (trailing zeros, exact decimal string, etc.).
    $textCells = @(
        "D5",
        "D6",
        "D9",
        "D10",
        "D14",
        "D19",
        "D21",
        "D23",
        "D24",
        "D25",
        "D27",
        "D30",
        "D32",
        "D33",
        "D34",
        "D35",
        "D38",
        "D39",
        "D40",
        "D42",
        "D43",
        "D44",
        "D50",
        "D51",
    )
    foreach ($addr in $textCells) {
        $ws.Range($addr).NumberFormat = "@"
    }

    $ws.Range("D2").Value = "51.375.90"
    $ws.Range("E2").Value = "  -1.02%  "
    $ws.Range("D3").Value = "2.770.89"
    $ws.Range("E3").Value = "  -0.62%  "
    $ws.Range("D5").Value = "352.12"
    $ws.Range("E5").Value = "  -1.87%  "
    $ws.Range("D6").Value = "108.05"
    $ws.Range("E6").Value = "  -1.51%  "
    $ws.Range("E7").Value = "  -1.53%  "
    $ws.Range("E8").Value = "  +0.08%  "
    $ws.Range("D9").Value = "0.618"
    $ws.Range("E9").Value = "  +4.85%  "
    $ws.Range("D10").Value = "39.21"
    $ws.Range("E10").Value = "  -1.79%  "
    $ws.Range("E11").Value = "  +1.70%  "
    $ws.Range("E12").Value = "  -1.44%  "
    $ws.Range("E13").Value = "  +1.75%  "
    $ws.Range("D14").Value = "7.79"
    $ws.Range("E14").Value = "  +3.19%  "
    $ws.Range("D15").Value = "3.209.68"
    $ws.Range("E15").Value = "  -0.45%  "
    $ws.Range("D16").Value = "2.798.48"
    $ws.Range("E16").Value = "  +0.40%  "
    $ws.Range("E17").Value = "  -1.51%  "
    $ws.Range("D18").Value = "51.418.72"
    $ws.Range("E18").Value = "  -0.88%  "
    $ws.Range("D19").Value = "7.71"
    $ws.Range("E19").Value = "  +3.04%  "
    $ws.Range("D21").Value = "13.32"
    $ws.Range("E21").Value = "  +1.59%  "
    $ws.Range("D22").Value = "0.0₃0967"
    $ws.Range("E22").Value = "  -0.81%  "
    $ws.Range("D23").Value = "70.39"
    $ws.Range("E23").Value = "  +0.17%  "
    $ws.Range("D24").Value = "266.16"
    $ws.Range("E24").Value = "  -1.41%  "
    $ws.Range("D25").Value = "2.77"
    $ws.Range("E25").Value = "  +0.99%  "
    $ws.Range("E26").Value = "  -0.17%  "
    $ws.Range("D27").Value = "25.84"
    $ws.Range("E27").Value = "  -2.61%  "
    $ws.Range("E28").Value = "  +0.64%  "
    $ws.Range("E29").Value = "  -0.09%  "
    $ws.Range("D30").Value = "36.92"
    $ws.Range("E30").Value = "  +8.58%  "
    $ws.Range("E31").Value = "  -2.27%  "
    $ws.Range("D32").Value = "6.18"
    $ws.Range("E32").Value = "  +7.93%  "
    $ws.Range("D33").Value = "52.00"
    $ws.Range("E33").Value = "  +0.09%  "
    $ws.Range("D34").Value = "0.0445"
    $ws.Range("E34").Value = "  -5.45%  "
    $ws.Range("D35").Value = "5.55"
    $ws.Range("E35").Value = "  +6.13%  "
    $ws.Range("E36").Value = "  -0.06%  "
    $ws.Range("E37").Value = "  -0.66%  "
    $ws.Range("D38").Value = "18.56"
    $ws.Range("E38").Value = "  -2.43%  "
    $ws.Range("D39").Value = "3.08"
    $ws.Range("E39").Value = "  -4.06%  "
    $ws.Range("D40").Value = "1.95"
    $ws.Range("E40").Value = "  -2.11%  "
    $ws.Range("E41").Value = "  -1.13%  "
    $ws.Range("D42").Value = "2.49"
    $ws.Range("E42").Value = "  -3.88%  "
    $ws.Range("D43").Value = "120.11"
    $ws.Range("E43").Value = "  +0.36%  "
    $ws.Range("D44").Value = "22.05"
    $ws.Range("E44").Value = "  +0.84%  "
    $ws.Range("E45").Value = "  -2.55%  "
    $ws.Range("D46").Value = "2.134.59"
    $ws.Range("E46").Value = "  +2.37%  "
    $ws.Range("E47").Value = "  +1.53%  "
    $ws.Range("E48").Value = "  +5.41%  "
    $ws.Range("E49").Value = "  +17.28%  "
    $ws.Range("D50").Value = "5.48"
    $ws.Range("E50").Value = "  -5.53%  "
    $ws.Range("D51").Value = "0.895"
    $ws.Range("E51").Value = "  -6.51%  "
